$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.191.25"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "'2.885.16"
$ws.Range("E3").Value = "  +4.22%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'352.26"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'112.01"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("D7").Value = "'0.559"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "'40.32"
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "'20.15"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("E16").Value = "  +8.31%  "
$ws.Range("D17").Value = "'2.893.13"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "'52.203.85"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +9.42%  "
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'13.75"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").Value = "'0.0₃0984"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "'71.27"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "'270.03"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "'26.44"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'0.164"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").Value = "'38.82"
$ws.Range("E30").Value = "  +5.45%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'6.33"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("E33").Value = "  +8.21%  "
$ws.Range("D34").Value = "'53.15"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").Value = "'0.0931"
$ws.Range("E35").Value = "  +11.37%  "
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  +8.02%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("D41").Value = "'2.63"
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("E42").Value = "  +2.97%  "
$ws.Range("D43").Value = "'22.43"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("D44").Value = "'121.71"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("E46").Value = "  +9.38%  "
$ws.Range("D47").Value = "'2.189.93"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "'2.47"
$ws.Range("D49").Value = "'0.263"
$ws.Range("E49").Value = "  +18.22%  "
$ws.Range("D50").Value = "'0.959"
$ws.Range("E50").Value = "  +7.16%  "
$ws.Range("D51").Value = "'0.0324"
$ws.Range("E51").Value = "  +11.16%  "
